# Edit: insert a new weekly record at row 11 (Fruta / hortaliza, semanal)
# This shifts all existing data rows 11-43 down to rows 12-44, and the new
# row 11 receives a new data point (date 45030, volume 300, prices 1500/1500/1500/1500).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 11, pushing existing rows 11-43 down to 12-44.
$ws.Rows(11).Insert()

# Populate the new row 11 with the new record's data.
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = 45030
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 100112044
$ws.Range("G11").Value = "Perejil"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = 1500
$ws.Range("N11").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 1500
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"

# Make sure the date cell keeps the same number format as the other dates.
$ws.Range("D11").NumberFormat = $ws.Range("D12").NumberFormat
